$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 27.97007004671944
$ws.Cells.Item(2, 3).Value = 31.68867800851428
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 28.03449121931777
$ws.Cells.Item(3, 3).Value = 31.69046039505187
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 28.01252619555745
$ws.Cells.Item(4, 3).Value = 31.65603535047414
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 27.97261962171509
$ws.Cells.Item(5, 3).Value = 31.67257253726252
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 27.99611430079139
$ws.Cells.Item(6, 3).Value = 31.71361694200338
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 28.04799599133151
$ws.Cells.Item(7, 3).Value = 31.71961988311895
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 27.99942927596777
$ws.Cells.Item(8, 3).Value = 31.7035024573028
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 28.01599650206708
$ws.Cells.Item(9, 3).Value = 31.70416323438609
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 27.98671635152283
$ws.Cells.Item(10, 3).Value = 31.69690310223311
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 27.98420113361988
$ws.Cells.Item(11, 3).Value = 31.677339656012
$ws.Cells.Item(12, 2).Value = 28.00201606386102
$ws.Cells.Item(12, 3).Value = 31.69228915663592

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 19.76205740453999
$ws.Cells.Item(2, 3).Value = 27.74693997578355
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 19.73664198286686
$ws.Cells.Item(3, 3).Value = 27.67832990789309
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 19.76095486155358
$ws.Cells.Item(4, 3).Value = 27.72475059421731
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 19.74063659989309
$ws.Cells.Item(5, 3).Value = 27.67851551070935
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 19.75833044360287
$ws.Cells.Item(6, 3).Value = 27.69141399864551
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19.76292998557473
$ws.Cells.Item(7, 3).Value = 27.72473321257943
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 19.75873829569051
$ws.Cells.Item(8, 3).Value = 27.72685943865804
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 19.7618823372115
$ws.Cells.Item(9, 3).Value = 27.67550658962473
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 19.74532170726462
$ws.Cells.Item(10, 3).Value = 27.71328803654262
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 19.767482925493
$ws.Cells.Item(11, 3).Value = 27.74172195830554
$ws.Cells.Item(12, 2).Value = 19.75549765436907
$ws.Cells.Item(12, 3).Value = 27.71020592229592

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 14.76077522084448
$ws.Cells.Item(2, 3).Value = 23.72092076782851
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 14.77326372306804
$ws.Cells.Item(3, 3).Value = 23.76438411778874
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 14.76042299534725
$ws.Cells.Item(4, 3).Value = 23.75829417070175
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 14.76488345217772
$ws.Cells.Item(5, 3).Value = 23.77154423494819
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 14.78040284802067
$ws.Cells.Item(6, 3).Value = 23.75036952552434
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 14.77004262157685
$ws.Cells.Item(7, 3).Value = 23.75436706937219
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 14.79261972188354
$ws.Cells.Item(8, 3).Value = 23.77086714609264
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 14.75452863854599
$ws.Cells.Item(9, 3).Value = 23.80531403503588
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 14.77067442157021
$ws.Cells.Item(10, 3).Value = 23.77261367942483
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 14.7860912511361
$ws.Cells.Item(11, 3).Value = 23.80259942904271
$ws.Cells.Item(12, 2).Value = 14.77137048941708
$ws.Cells.Item(12, 3).Value = 23.76712741757598

Write-Host "Update complete"
